$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 6 and 7 (Khorfakkan-Hatta / Al Wasl-Emirates Club) were reordered.
#    Columns A-E (index/pais/torneio/temporada/data_partida) stay the same;
#    columns F-V (match details) swap between the two rows.
# ---------------------------------------------------------------------------
$row6 = @($ws.Cells.Item(6,6).Value2, $ws.Cells.Item(6,7).Value2, $ws.Cells.Item(6,8).Value2, $ws.Cells.Item(6,9).Value2, $ws.Cells.Item(6,10).Value2, $ws.Cells.Item(6,11).Value2, $ws.Cells.Item(6,12).Value2, $ws.Cells.Item(6,13).Value2, $ws.Cells.Item(6,14).Value2, $ws.Cells.Item(6,15).Value2, $ws.Cells.Item(6,16).Value2, $ws.Cells.Item(6,17).Value2, $ws.Cells.Item(6,18).Value2, $ws.Cells.Item(6,19).Value2, $ws.Cells.Item(6,20).Value2, $ws.Cells.Item(6,21).Value2, $ws.Cells.Item(6,22).Value2)

$row7 = @($ws.Cells.Item(7,6).Value2, $ws.Cells.Item(7,7).Value2, $ws.Cells.Item(7,8).Value2, $ws.Cells.Item(7,9).Value2, $ws.Cells.Item(7,10).Value2, $ws.Cells.Item(7,11).Value2, $ws.Cells.Item(7,12).Value2, $ws.Cells.Item(7,13).Value2, $ws.Cells.Item(7,14).Value2, $ws.Cells.Item(7,15).Value2, $ws.Cells.Item(7,16).Value2, $ws.Cells.Item(7,17).Value2, $ws.Cells.Item(7,18).Value2, $ws.Cells.Item(7,19).Value2, $ws.Cells.Item(7,20).Value2, $ws.Cells.Item(7,21).Value2, $ws.Cells.Item(7,22).Value2)

for ($i = 0; $i -lt 17; $i++) {
    $ws.Cells.Item(6, 6 + $i).Value = $row7[$i]
    $ws.Cells.Item(7, 6 + $i).Value = $row6[$i]
}

# ---------------------------------------------------------------------------
# 2) Append new row 48 (Ittihad Kalba vs Shabab Al-Ahli Dubai), matching the
#    formatting already used by the existing data rows (bold/bordered index
#    cell in column A, datetime format in column E).
# ---------------------------------------------------------------------------
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("E47").Copy()
$ws.Range("E48").PasteSpecial(-4122)

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "united-arab-emirates"
$ws.Range("C48").Value = "uae-league"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45233.6875
$ws.Range("F48").Value = "Ittihad Kalba"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "Shabab Al-Ahli Dubai"
$ws.Range("I48").Value = 3
$ws.Range("J48").Value = 3.9
$ws.Range("K48").Value = "30/10/2023 18:42"
$ws.Range("L48").Value = 4.27
$ws.Range("M48").Value = "03/11/2023 16:28"
$ws.Range("N48").Value = 4
$ws.Range("O48").Value = "30/10/2023 18:42"
$ws.Range("P48").Value = 4.44
$ws.Range("Q48").Value = "03/11/2023 16:24"
$ws.Range("R48").Value = 1.82
$ws.Range("S48").Value = "30/10/2023 18:42"
$ws.Range("T48").Value = 1.69
$ws.Range("U48").Value = "03/11/2023 16:24"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ittihad-kalba-shabab-al-ahli-dubai/fqBKarJj/"
